$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 45457956
$ws.Range("I98").Value = 58824804
$ws.Range("J98").Value = 10679
$ws.Range("K98").Value = 58824804
$ws.Range("L98").Value = 10679
$ws.Range("M98").Value = -58823306
$ws.Range("N98").Value = -13675

$ws.Range("H105").Value = 88944.55
$ws.Range("J105").Value = 88944.55
$ws.Range("L105").Value = 88944.55
$ws.Range("N105").Value = -95932.55

$ws.Range("H106").Value = 19873.273
$ws.Range("I106").Value = 2066.111
$ws.Range("K106").Value = 2066.111
$ws.Range("M106").Value = -1435.111

$ws.Range("H116").Value = 5698.1177
$ws.Range("I116").Value = 5237.4
$ws.Range("K116").Value = 5237.4
$ws.Range("M116").Value = -1795.4

$ws.Range("H122").Value = 45457956
$ws.Range("I122").Value = 58824804
$ws.Range("J122").Value = 10679
$ws.Range("K122").Value = 176474412
$ws.Range("L122").Value = 32037
$ws.Range("M122").Value = -176471962
$ws.Range("N122").Value = -36937

$ws.Range("H137").Value = 9189.416999999999
$ws.Range("I137").Value = 2995.375
$ws.Range("K137").Value = 8986.125
$ws.Range("M137").Value = -6436.125

$ws.Range("H138").Value = 4507.959
$ws.Range("I138").Value = 3924.7646
$ws.Range("J138").Value = 4817.7812
$ws.Range("K138").Value = 11774.2938
$ws.Range("L138").Value = 14453.3436
$ws.Range("M138").Value = -6634.293799999999
$ws.Range("N138").Value = -24733.3436

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 800
$ws.Range("I6").Value = 800
$ws.Range("K6").Value = 800
$ws.Range("M6").Value = -627

$ws.Range("H101").Value = 85670.5
$ws.Range("J101").Value = 85670.5
$ws.Range("L101").Value = 85670.5
$ws.Range("N101").Value = -92160.5

$ws.Range("H102").Value = 7998.5
$ws.Range("I102").Value = 8213.6
$ws.Range("K102").Value = 8213.6
$ws.Range("M102").Value = -6591.6

$ws.Range("H132").Value = 6323.0176
$ws.Range("I132").Value = 3663.4285
$ws.Range("K132").Value = 10990.2855
$ws.Range("M132").Value = -8460.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 95780
$ws.Range("J50").Value = 95780
$ws.Range("L50").Value = 95780
$ws.Range("N50").Value = -96928

$ws.Range("H51").Value = 299775
$ws.Range("J51").Value = 299775
$ws.Range("L51").Value = 299775
$ws.Range("N51").Value = -300757

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 46529.668
$ws.Range("J18").Value = 46529.668
$ws.Range("L18").Value = 46529.668
$ws.Range("N18").Value = -46989.668

$ws.Range("H28").Value = 30918.223
$ws.Range("J28").Value = 30918.223
$ws.Range("L28").Value = 30918.223
$ws.Range("N28").Value = -31408.223

$ws.Range("H31").Value = 941113.5600000001
$ws.Range("I31").Value = 13653.077
$ws.Range("J31").Value = 1945862.4
$ws.Range("K31").Value = 13653.077
$ws.Range("L31").Value = 1945862.4
$ws.Range("M31").Value = -13358.077
$ws.Range("N31").Value = -1946452.4

$ws.Range("H34").Value = 941113.5600000001
$ws.Range("I34").Value = 13653.077
$ws.Range("J34").Value = 1945862.4
$ws.Range("K34").Value = 13653.077
$ws.Range("L34").Value = 1945862.4
$ws.Range("M34").Value = -13451.077
$ws.Range("N34").Value = -1946266.4

$ws.Range("H58").Value = 1998.5264
$ws.Range("I58").Value = 1548.5
$ws.Range("J58").Value = 2498.5557
$ws.Range("K58").Value = 1548.5
$ws.Range("L58").Value = 2498.5557
$ws.Range("M58").Value = -1345.5
$ws.Range("N58").Value = -2904.5557

$ws.Range("H62").Value = 3217.25
$ws.Range("I62").Value = 3031.3333
$ws.Range("K62").Value = 3031.3333
$ws.Range("M62").Value = -2407.3333

$ws.Range("H65").Value = 3217.25
$ws.Range("I65").Value = 3031.3333
$ws.Range("K65").Value = 15156.6665
$ws.Range("M65").Value = -12036.6665

$ws.Range("H107").Value = 772.9524
$ws.Range("I107").Value = 753.7895
$ws.Range("J107").Value = 955
$ws.Range("K107").Value = 753.7895
$ws.Range("L107").Value = 955
$ws.Range("M107").Value = 1166.2105
$ws.Range("N107").Value = -4795

$ws.Range("H122").Value = 2425.4783
$ws.Range("I122").Value = 1547.4375
$ws.Range("K122").Value = 4642.3125
$ws.Range("M122").Value = -2192.3125

$ws.Range("H132").Value = 3309.1143
$ws.Range("J132").Value = 5057.2
$ws.Range("L132").Value = 15171.6
$ws.Range("N132").Value = -20231.6

$ws.Range("H136").Value = 1998.5264
$ws.Range("I136").Value = 1548.5
$ws.Range("J136").Value = 2498.5557
$ws.Range("K136").Value = 4645.5
$ws.Range("L136").Value = 7495.6671
$ws.Range("M136").Value = -2095.5
$ws.Range("N136").Value = -12595.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7408.8096
$ws.Range("I56").Value = 7408.8096
$ws.Range("K56").Value = 7408.8096
$ws.Range("M56").Value = -6878.8096

$ws.Range("H68").Value = 1942.7812
$ws.Range("J68").Value = 1930.5454
$ws.Range("L68").Value = 5791.6362
$ws.Range("N68").Value = -7413.6362

$ws.Range("H71").Value = 1942.7812
$ws.Range("J71").Value = 1930.5454
$ws.Range("L71").Value = 17374.9086
$ws.Range("N71").Value = -25486.9086

$ws.Range("H107").Value = 551.7143
$ws.Range("I107").Value = 551.7143
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1655.1429
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 264.8571000000002
$ws.Range("N107").ClearContents()

$ws.Range("H134").Value = 25004896
$ws.Range("J134").Value = 8000
$ws.Range("L134").Value = 24000
$ws.Range("N134").Value = -34140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 36651.332
$ws.Range("I5").Value = 34974.5
$ws.Range("J5").Value = 40005
$ws.Range("K5").Value = 34974.5
$ws.Range("L5").Value = 40005
$ws.Range("M5").Value = -34862.5
$ws.Range("N5").Value = -40229

$ws.Range("H106").Value = 113940
$ws.Range("J106").Value = 113940
$ws.Range("L106").Value = 113940
$ws.Range("N106").Value = -116464

$ws.Range("H122").Value = 5337.409
$ws.Range("I122").Value = 5232.3687
$ws.Range("K122").Value = 15697.1061
$ws.Range("M122").Value = -13247.1061

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 74613.92999999999
$ws.Range("I7").Value = 1857.5
$ws.Range("J7").Value = 256505
$ws.Range("K7").Value = 1857.5
$ws.Range("L7").Value = 256505
$ws.Range("M7").Value = -1745.5
$ws.Range("N7").Value = -256729

$ws.Range("H22").Value = 3907.1667
$ws.Range("I22").Value = 4842.154
$ws.Range("J22").Value = 3192.1765
$ws.Range("K22").Value = 4842.154
$ws.Range("L22").Value = 3192.1765
$ws.Range("M22").Value = -4547.154
$ws.Range("N22").Value = -3782.1765

$ws.Range("H27").Value = 3907.1667
$ws.Range("I27").Value = 4842.154
$ws.Range("J27").Value = 3192.1765
$ws.Range("K27").Value = 4842.154
$ws.Range("L27").Value = 3192.1765
$ws.Range("M27").Value = -4735.154
$ws.Range("N27").Value = -3406.1765

$ws.Range("J40").Value = 4999.8887
$ws.Range("L40").Value = 4999.8887
$ws.Range("N40").Value = -5271.8887

$ws.Range("H46").Value = 2294.0557
$ws.Range("I46").Value = 2249.3
$ws.Range("K46").Value = 2249.3
$ws.Range("M46").Value = -2061.3

$ws.Range("H61").Value = 1585.8889
$ws.Range("I61").Value = 1513.2916
$ws.Range("K61").Value = 1513.2916
$ws.Range("M61").Value = -1311.2916

$ws.Range("H93").Value = 45456756
$ws.Range("I93").Value = 58825156
$ws.Range("K93").Value = 58825156
$ws.Range("M93").Value = -58823908

$ws.Range("H103").Value = 32714.334
$ws.Range("J103").Value = 32714.334
$ws.Range("L103").Value = 32714.334
$ws.Range("N103").Value = -35058.334

$ws.Range("H106").Value = 32880
$ws.Range("J106").Value = 32880
$ws.Range("L106").Value = 32880
$ws.Range("N106").Value = -35404

$ws.Range("H113").Value = 1585.8889
$ws.Range("I113").Value = 1513.2916
$ws.Range("K113").Value = 1513.2916
$ws.Range("M113").Value = 656.7084

$ws.Range("H126").Value = 74613.92999999999
$ws.Range("I126").Value = 1857.5
$ws.Range("J126").Value = 256505
$ws.Range("K126").Value = 5572.5
$ws.Range("L126").Value = 769515
$ws.Range("M126").Value = -3102.5
$ws.Range("N126").Value = -774455

$ws.Range("H136").Value = 73092.45
$ws.Range("I136").Value = 66608.25
$ws.Range("K136").Value = 199824.75
$ws.Range("M136").Value = -197274.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 24500
$ws.Range("I32").Value = 24500
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 24500
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -24183
$ws.Range("N32").ClearContents()

$ws.Range("H34").Value = 29999
$ws.Range("J34").Value = 29999
$ws.Range("L34").Value = 29999
$ws.Range("N34").Value = -30405

$ws.Range("H104").Value = 81396.664
$ws.Range("J104").Value = 81396.664
$ws.Range("L104").Value = 81396.664
$ws.Range("N104").Value = -88384.664
